$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: petr4 - updated metrics ---
$ws.Range("B2").Value = 37.34
$ws.Range("D2").Value = 71.98322222222224
$ws.Range("E2").Value = 77.034894690653
$ws.Range("F2").Value = 22.45
$ws.Range("G2").Value = 8.44
$ws.Range("H2").Value = 17.35
$ws.Range("I2").Value = 22.60310658810926
$ws.Range("J2").Value = 83.69041242635244
$ws.Range("K2").Value = 4.43
$ws.Range("L2").Value = 1.2
$ws.Range("M2").Value = 27.02
$ws.Range("N2").Value = 17.93
$ws.Range("P2").Value = 3.07
$ws.Range("R2").Value = 0

# --- Row 3: vale3 replaced by cmig3 ---
$ws.Range("A3").Value = "cmig3"
$ws.Range("B3").Value = 12.93
$ws.Range("C3").Value = 8.880000000000001
$ws.Range("D3").Value = 15.07063333333333
$ws.Range("E3").Value = 19.63705680594727
$ws.Range("F3").Value = 14.81
$ws.Range("G3").Value = 1.93
$ws.Range("H3").Value = 10.49
$ws.Range("I3").Value = 14.92652745552978
$ws.Range("J3").Value = 68.67749419953597
$ws.Range("K3").Value = 6.71
$ws.Range("L3").Value = 1.46
$ws.Range("M3").Value = 21.72
$ws.Range("N3").Value = 16.13
$ws.Range("O3").Value = 80
$ws.Range("P3").Value = 4.76
$ws.Range("Q3").Value = 0.28
$ws.Range("R3").Value = 0.86
$ws.Range("S3").Value = 1.15

# --- Row 4: cmig3 replaced by sbsp3 ---
$ws.Range("A4").Value = "sbsp3"
$ws.Range("B4").Value = 88.65000000000001
$ws.Range("C4").Value = 44.89
$ws.Range("D4").Value = 16.0555
$ws.Range("E4").Value = 72.95773947978377
$ws.Range("F4").Value = 13.62
$ws.Range("G4").Value = 5.27
$ws.Range("H4").Value = 1.63
$ws.Range("I4").Value = 5.944726452340665
$ws.Range("J4").Value = 50.63733784545968
$ws.Range("K4").Value = 16.82
$ws.Range("L4").Value = 1.97
$ws.Range("M4").Value = 11.73
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 100
$ws.Range("P4").Value = 7.94
$ws.Range("Q4").Value = 0.52
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 1.23

# --- Shrink conditional-formatting ranges from row 5 down to row 4 ---
$cols = @("F", "Q", "N", "R", "O", "K", "L", "S", "M", "P", "I")
foreach ($col in $cols) {
    $rangeOld = $col + "2:" + $col + "5"
    $rangeNew = $col + "2:" + $col + "4"
    $fcs = $ws.Range($rangeOld).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($rangeNew))
    }
}

# --- Row 5: azul4 removed entirely ---
$ws.Rows("5").Delete()
